$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/373537613e69257f2c6ac12ff3090b1f69ab2d3f/e2e/5c5d57c7-208e-47aa-957f-7aa028c60097.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bffbf4dfa822bda499fd971073ac409d187e8b23/e2e/5c5d57c7-208e-47aa-957f-7aa028c60097.md."

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen columns I, J, P to 40
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40
$wsZh.Columns.Item(16).ColumnWidth = 40

# Fill in row 6 handback-report details
$wsZh.Range("I6").Value = "5c5d57c7-208e-47aa-957f-7aa028c60097.md"
$wsZh.Range("J6").Value = "5c5d57c7-208e-47aa-957f-7aa028c60097.10c5ae0019643c43d770f41225bcd6b65b91a091.zh-cn.xlf"
$wsZh.Range("K6").Value = "2016-09-06 10:05:38"
$wsZh.Range("P6").Value = $errorDetail

# Hyperlink on I6 pointing at the handback markdown file (same target as A6),
# styled to match the workbook's existing hyperlink look (underline + blue)
$wsZh.Hyperlinks.Add($wsZh.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bffbf4dfa822bda499fd971073ac409d187e8b23/e2e/5c5d57c7-208e-47aa-957f-7aa028c60097.md", "", "", "5c5d57c7-208e-47aa-957f-7aa028c60097.md")
$wsZh.Range("I6").Font.Underline = $true
$wsZh.Range("I6").Font.Color = 15570276

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")

# Widen columns I, J, P to 40
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40
$wsDe.Columns.Item(16).ColumnWidth = 40

# Fill in row 6 handback-report details
$wsDe.Range("I6").Value = "5c5d57c7-208e-47aa-957f-7aa028c60097.md"
$wsDe.Range("J6").Value = "5c5d57c7-208e-47aa-957f-7aa028c60097.10c5ae0019643c43d770f41225bcd6b65b91a091.de-de.xlf"
$wsDe.Range("K6").Value = "2016-09-06 10:05:55"
$wsDe.Range("P6").Value = $errorDetail

# Hyperlink on I6 pointing at the handback markdown file (same target as A6),
# styled to match the workbook's existing hyperlink look (underline + blue)
$wsDe.Hyperlinks.Add($wsDe.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bffbf4dfa822bda499fd971073ac409d187e8b23/e2e/5c5d57c7-208e-47aa-957f-7aa028c60097.md", "", "", "5c5d57c7-208e-47aa-957f-7aa028c60097.md")
$wsDe.Range("I6").Font.Underline = $true
$wsDe.Range("I6").Font.Color = 15570276
